$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simulation score")
$ws.Name = "Performance score"
